$wb = $excel.ActiveWorkbook

# Target "want-to-go" counts (column F) for the rows that remain after the
# 2024-04-21 "合肥·银魂only" entry (previously row 2) is removed and every
# later row shifts up by one. Keyed by the NEW row number (2-based data rows).
$targetF = @{
    2  = 106
    3  = 398
    4  = 11670
    5  = 853
    6  = 116
    7  = 16
    8  = 84
    9  = 148
    10 = 171
    11 = 24
    12 = 50
    13 = 53
    14 = 129
    15 = 34
    16 = 336
    17 = 1370
    18 = 78
    19 = 904
    20 = 111
}

# Sheets that contain the "展览" style listing with the expired 2024-04-21 row.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

    # Remove the expired "2024-04-21" event (row 2); this shifts every
    # subsequent row up by one automatically.
    $ws.Rows.Item(2).Delete()

    $newLastRow = $lastRow - 1

    # Excel's row delete also shifted the plain numeric values that were
    # sitting in column A (the running index 0,1,2,...), so renumber them
    # back into a clean sequential series starting at 0.
    for ($r = 1; $r -le $newLastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the "want to go" counters (column F) for the rows that moved
    # up, matching the latest scrape values.
    foreach ($rowKey in $targetF.Keys) {
        if ([int]$rowKey -le $newLastRow) {
            $ws.Cells.Item([int]$rowKey, 6).Value = $targetF[$rowKey]
        }
    }
}
